$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell T1 (continues the 0..17 numbering with 18) ---
$ws.Range('T1').Value = 18

# --- New row labels, entered first so they land in the shared-string table
#     in the same relative order the authored workbook uses (OffsetF,
#     OffsetA, RD Single, TD Single, then 1Pair-B). ---
$ws.Range('B7').Value = 'OffsetF'
$ws.Range('B8').Value = 'OffsetA'
$ws.Range('B9').Value = 'RD Single'
$ws.Range('B10').Value = 'TD Single'

# --- Row 2 (HKL reflection labels): insert "1Pair-B" before the existing
#     "2Pairs-A" group, shifting the remaining labels right by one column,
#     and add the trailing "MaxUnique" label in the new T column. ---
$ws.Range('L2').Value = '1Pair-B'
$ws.Range('M2').Value = '2Pairs-A'
$ws.Range('N2').Value = '2Pairs-B'
$ws.Range('O2').Value = '3Pairs-A'
$ws.Range('P2').Value = '3Pairs-B'
$ws.Range('Q2').Value = '3Pairs-C'
$ws.Range('R2').Value = '4Pairs'
$ws.Range('S2').Value = '5A4F'
$ws.Range('T2').Value = 'MaxUnique'

# --- Row 3 (Equal Angle): refreshed intensities + new T column ---
$ws.Range('C3').Value = 0.9153170028818444
$ws.Range('D3').Value = 1.251433717579251
$ws.Range('E3').Value = 0.9428458213256484
$ws.Range('F3').Value = 0.9153170028818444
$ws.Range('G3').Value = 1.132118155619597
$ws.Range('H3').Value = 0.8583933717579251
$ws.Range('I3').Value = 0.9364337175792508
$ws.Range('J3').Value = 1.251433717579251
$ws.Range('K3').Value = 0.9153170028818444
$ws.Range('L3').Value = 0.9428458213256484
$ws.Range('M3').Value = 1.09713976945245
$ws.Range('N3').Value = 1.09713976945245
$ws.Range('O3').Value = 1.108799231508165
$ws.Range('P3').Value = 1.036532180595581
$ws.Range('Q3').Value = 1.036532180595581
$ws.Range('R3').Value = 1.006228386167147
$ws.Range('S3').Value = 1.006228386167147
$ws.Range('T3').Value = 1.006090297790586

# --- Row 4 (CLR): refreshed intensities + new T column ---
$ws.Range('C4').Value = 0.9967415414473634
$ws.Range('D4').Value = 1.014824979741739
$ws.Range('E4').Value = 0.989674153860961
$ws.Range('F4').Value = 0.9967415414473634
$ws.Range('G4').Value = 1.010550431152634
$ws.Range('H4').Value = 0.9755317116315174
$ws.Range('I4').Value = 0.9907461123635702
$ws.Range('J4').Value = 1.014824979741739
$ws.Range('K4').Value = 0.9967415414473634
$ws.Range('L4').Value = 0.989674153860961
$ws.Range('M4').Value = 1.00224956680135
$ws.Range('N4').Value = 1.00224956680135
$ws.Range('O4').Value = 1.005016521585111
$ws.Range('P4').Value = 1.000413558350021
$ws.Range('Q4').Value = 1.000413558350021
$ws.Range('R4').Value = 0.9994955541243569
$ws.Range('S4').Value = 0.9994955541243569
$ws.Range('T4').Value = 0.996344821699631

# --- Row 5 (BT8Hex): refreshed intensities + new T column ---
$ws.Range('C5').Value = 0.9920719306253711
$ws.Range('D5').Value = 1.036616406673343
$ws.Range('E5').Value = 0.9841311122394232
$ws.Range('F5').Value = 0.9920719306253711
$ws.Range('G5').Value = 1.022569247517488
$ws.Range('H5').Value = 0.9597013650732229
$ws.Range('I5').Value = 0.9851706286313915
$ws.Range('J5').Value = 1.036616406673343
$ws.Range('K5').Value = 0.9920719306253711
$ws.Range('L5').Value = 0.9841311122394232
$ws.Range('M5').Value = 1.010373759456383
$ws.Range('N5').Value = 1.010373759456383
$ws.Range('O5').Value = 1.014438922143418
$ws.Range('P5').Value = 1.004273149846046
$ws.Range('Q5').Value = 1.004273149846046
$ws.Range('R5').Value = 1.001222845040877
$ws.Range('S5').Value = 1.001222845040877
$ws.Range('T5').Value = 0.9967101151267066

# --- Row 6 (Spiral): refreshed intensities + new T column ---
$ws.Range('C6').Value = 0.9980043998191661
$ws.Range('D6').Value = 0.9985180596447223
$ws.Range('E6').Value = 0.9950840607171098
$ws.Range('F6').Value = 0.9980043998191661
$ws.Range('G6').Value = 1.0002790861062
$ws.Range('H6').Value = 0.9915136988992173
$ws.Range('I6').Value = 0.995149242567458
$ws.Range('J6').Value = 0.9985180596447223
$ws.Range('K6').Value = 0.9980043998191661
$ws.Range('L6').Value = 0.9950840607171098
$ws.Range('M6').Value = 0.9968010601809161
$ws.Range('N6').Value = 0.9968010601809161
$ws.Range('O6').Value = 0.9979604021560107
$ws.Range('P6').Value = 0.9972021733936661
$ws.Range('Q6').Value = 0.9972021733936661
$ws.Range('R6').Value = 0.9974027300000411
$ws.Range('S6').Value = 0.9974027300000411
$ws.Range('T6').Value = 0.9964247579589789

# --- Row 7: label now resolves to "OffsetF" (new category inserted into the
#     shared-string table ahead of it; see $ws.Range('B7') above) with
#     refreshed intensities + new T column ---
$ws.Range('C7').Value = 1.169935231624025
$ws.Range('D7').Value = 0.501275923595384
$ws.Range('E7').Value = 1.088452508488073
$ws.Range('F7').Value = 1.169935231624025
$ws.Range('G7').Value = 0.736759667576416
$ws.Range('H7').Value = 1.228842748610638
$ws.Range('I7').Value = 1.119103383396861
$ws.Range('J7').Value = 0.501275923595384
$ws.Range('K7').Value = 1.169935231624025
$ws.Range('L7').Value = 1.088452508488073
$ws.Range('M7').Value = 0.7948642160417285
$ws.Range('N7').Value = 0.7948642160417285
$ws.Range('O7').Value = 0.7754960332199577
$ws.Range('P7').Value = 0.9198878879024939
$ws.Range('Q7').Value = 0.9198878879024939
$ws.Range('R7').Value = 0.9823997238328765
$ws.Range('S7').Value = 0.9823997238328765
$ws.Range('T7').Value = 0.9740615772152328

# --- New row 8: OffsetA (label set above) ---
$ws.Range('A8').Value = 6
$ws.Range('C8').Value = 0.989119581040478
$ws.Range('D8').Value = 0.9088275117754978
$ws.Range('E8').Value = 1.030035410330014
$ws.Range('F8').Value = 0.989119581040478
$ws.Range('G8').Value = 0.9359611356219023
$ws.Range('H8').Value = 1.114173862554452
$ws.Range('I8').Value = 1.024696613891124
$ws.Range('J8').Value = 0.9088275117754978
$ws.Range('K8').Value = 0.989119581040478
$ws.Range('L8').Value = 1.030035410330014
$ws.Range('M8').Value = 0.9694314610527561
$ws.Range('N8').Value = 0.9694314610527561
$ws.Range('O8').Value = 0.9582746859091382
$ws.Range('P8').Value = 0.9759941677153301
$ws.Range('Q8').Value = 0.9759941677153301
$ws.Range('R8').Value = 0.9792755210466171
$ws.Range('S8').Value = 0.9792755210466171
$ws.Range('T8').Value = 1.000469019202245

# --- New row 9: RD Single (label set above) ---
$ws.Range('A9').Value = 7
$ws.Range('C9').Value = 1.97
$ws.Range('D9').Value = 0.22
$ws.Range('E9').Value = 0.83
$ws.Range('F9').Value = 1.97
$ws.Range('G9').Value = 0.63
$ws.Range('H9').Value = 0.69
$ws.Range('I9').Value = 1.14
$ws.Range('J9').Value = 0.22
$ws.Range('K9').Value = 1.97
$ws.Range('L9').Value = 0.83
$ws.Range('M9').Value = 0.525
$ws.Range('N9').Value = 0.525
$ws.Range('O9').Value = 0.56
$ws.Range('P9').Value = 1.006666666666667
$ws.Range('Q9').Value = 1.006666666666667
$ws.Range('R9').Value = 1.2475
$ws.Range('S9').Value = 1.2475
$ws.Range('T9').Value = 0.9133333333333332

# --- New row 10: TD Single (label set above) ---
$ws.Range('A10').Value = 8
$ws.Range('C10').Value = 1.33
$ws.Range('D10').Value = 0.16
$ws.Range('E10').Value = 1.16
$ws.Range('F10').Value = 1.33
$ws.Range('G10').Value = 0.46
$ws.Range('H10').Value = 1.54
$ws.Range('I10').Value = 1.22
$ws.Range('J10').Value = 0.16
$ws.Range('K10').Value = 1.33
$ws.Range('L10').Value = 1.16
$ws.Range('M10').Value = 0.6599999999999999
$ws.Range('N10').Value = 0.6599999999999999
$ws.Range('O10').Value = 0.5933333333333333
$ws.Range('P10').Value = 0.8833333333333333
$ws.Range('Q10').Value = 0.8833333333333333
$ws.Range('R10').Value = 0.995
$ws.Range('S10').Value = 0.995
$ws.Range('T10').Value = 0.9783333333333334

# --- New row 11: HexGrid-90degTilt5degRes ---
$ws.Range('A11').Value = 9
$ws.Range('B11').Value = 'HexGrid-90degTilt5degRes'
$ws.Range('C11').Value = 0.9995825844195323
$ws.Range('D11').Value = 0.9947859823906227
$ws.Range('E11').Value = 0.9956587129011232
$ws.Range('F11').Value = 0.9995825844195323
$ws.Range('G11').Value = 0.9986076548371395
$ws.Range('H11').Value = 0.9926008551345861
$ws.Range('I11').Value = 0.9959131770519947
$ws.Range('J11').Value = 0.9947859823906227
$ws.Range('K11').Value = 0.9995825844195323
$ws.Range('L11').Value = 0.9956587129011232
$ws.Range('M11').Value = 0.9952223476458729
$ws.Range('N11').Value = 0.9952223476458729
$ws.Range('O11').Value = 0.9963507833762951
$ws.Range('P11').Value = 0.9966757599037593
$ws.Range('Q11').Value = 0.9966757599037593
$ws.Range('R11').Value = 0.9974024660327026
$ws.Range('S11').Value = 0.9974024660327026
$ws.Range('T11').Value = 0.9961914944558329

# --- Formatting: mirror the existing bold/centered/bordered header style
#     onto the newly-created header cell and the new row-index cells in
#     column A, the same way the rest of row 1 / column A is styled. ---
$ws.Range('S1').Copy()
$ws.Range('T1').PasteSpecial(-4122)

$ws.Range('A7').Copy()
$ws.Range('A8:A11').PasteSpecial(-4122)

$excel.CutCopyMode = $false
